$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 15, essentially a highlighted copy of row 12 with updated values
$ws.Range("A15:L15").Interior.Color = 65535

$ws.Range("A15").Value = $ws.Range("A12").Value2
$ws.Range("D15").Value = 3000
$ws.Range("E15").Value = 200
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 0.8
$ws.Range("H15").Value = 4000
$ws.Range("I15").Value = $ws.Range("I12").Value2
$ws.Range("J15").Value = $ws.Range("J12").Value2
$ws.Range("K15").Value = 112
$ws.Range("L15").Value = "pretty good"

$ws.Range("J19").Select()
